$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBE = New-Object 'object[,]' 24,4
$arrBE[0,0] = 15.5481441735457
$arrBE[0,1] = 12.3976867857438
$arrBE[0,2] = 5.96792099260609
$arrBE[0,3] = 16.53272111060208
$arrBE[1,0] = 14.80732950842882
$arrBE[1,1] = 11.59619543570442
$arrBE[1,2] = 5.847133426156652
$arrBE[1,3] = 15.58751913712275
$arrBE[2,0] = 14.33664346255196
$arrBE[2,1] = 11.09675394752501
$arrBE[2,2] = 5.773691212322553
$arrBE[2,3] = 14.98279392996644
$arrBE[3,0] = 14.14113299146339
$arrBE[3,1] = 10.88702080934155
$arrBE[3,2] = 5.743989483833317
$arrBE[3,3] = 14.73050559302048
$arrBE[4,0] = 14.10845392606804
$arrBE[4,1] = 10.85182376616609
$arrBE[4,2] = 5.739072539681684
$arrBE[4,3] = 14.688268433947
$arrBE[5,0] = 14.33402133382075
$arrBE[5,1] = 11.09395037415437
$arrBE[5,2] = 5.773289667308506
$arrBE[5,3] = 14.97941480395135
$arrBE[6,0] = 15.29614892130588
$arrBE[6,1] = 12.1172569150298
$arrBE[6,2] = 5.92615054732452
$arrBE[6,3] = 16.21201878222172
$arrBE[7,0] = 17.04713291467398
$arrBE[7,1] = 14.02004347292009
$arrBE[7,2] = 6.22968424711997
$arrBE[7,3] = 18.52819870552906
$arrBE[8,0] = 18.2395788168513
$arrBE[8,1] = 15.267591537137
$arrBE[8,2] = 6.452445787074532
$arrBE[8,3] = 20.19168851745649
$arrBE[9,0] = 18.75985759220933
$arrBE[9,1] = 15.80281928184702
$arrBE[9,2] = 6.553223165469839
$arrBE[9,3] = 20.90669558127356
$arrBE[10,0] = 18.95356902615478
$arrBE[10,1] = 16.00088056273671
$arrBE[10,2] = 6.591264846552171
$arrBE[10,3] = 21.17149113522515
$arrBE[11,0] = 18.91199854902262
$arrBE[11,1] = 15.95842950195859
$arrBE[11,2] = 6.583077868100959
$arrBE[11,3] = 21.11472724421126
$arrBE[12,0] = 18.7758612296765
$arrBE[12,1] = 15.81920630722452
$arrBE[12,2] = 6.556355515839001
$arrBE[12,3] = 20.9285997454808
$arrBE[13,0] = 18.69203937039664
$arrBE[13,1] = 15.73332744619367
$arrBE[13,2] = 6.539970449661986
$arrBE[13,3] = 20.81381594944768
$arrBE[14,0] = 18.20512061493702
$arrBE[14,1] = 15.23196557500471
$arrBE[14,2] = 6.445845026493198
$arrBE[14,3] = 20.14412477105941
$arrBE[15,0] = 17.90064202503795
$arrBE[15,1] = 14.91614443265502
$arrBE[15,2] = 6.38792986313771
$arrBE[15,3] = 19.72263150415201
$arrBE[16,0] = 17.72343282010558
$arrBE[16,1] = 14.73145184836543
$arrBE[16,2] = 6.354568351588591
$arrBE[16,3] = 19.4762703394299
$arrBE[17,0] = 17.66307934778457
$arrBE[17,1] = 14.66839511816029
$arrBE[17,2] = 6.343265351464196
$arrBE[17,3] = 19.39218065942177
$arrBE[18,0] = 17.93327061019536
$arrBE[18,1] = 14.95007864241637
$arrBE[18,2] = 6.394100529793449
$arrBE[18,3] = 19.76790658561762
$arrBE[19,0] = 18.81593865165496
$arrBE[19,1] = 15.86022466111552
$arrBE[19,2] = 6.564208099545667
$arrBE[19,3] = 20.9834314146336
$arrBE[20,0] = 19.37349723040085
$arrBE[20,1] = 16.42815299098328
$arrBE[20,2] = 6.674663627812241
$arrBE[20,3] = 21.74311279948103
$arrBE[21,0] = 19.07771918499286
$arrBE[21,1] = 16.12749243159752
$arrBE[21,2] = 6.615790007043737
$arrBE[21,3] = 21.34082236245071
$arrBE[22,0] = 17.91852594681535
$arrBE[22,1] = 14.93474671492008
$arrBE[22,2] = 6.391310971751499
$arrBE[22,3] = 19.74745030493806
$arrBE[23,0] = 16.58926006646703
$arrBE[23,1] = 13.53182415070685
$arrBE[23,2] = 6.147434359897931
$arrBE[23,3] = 17.8783551255268
$ws.Range("B2:E25").Value = $arrBE

$arrGI = New-Object 'object[,]' 24,3
$arrGI[0,0] = 34.71726866766956
$arrGI[0,1] = 14.60081989217757
$arrGI[0,2] = 19.95002442975634
$arrGI[1,0] = 33.88729454128996
$arrGI[1,1] = 14.56532119140542
$arrGI[1,2] = 19.94029781321222
$arrGI[2,0] = 33.38433052538839
$arrGI[2,1] = 14.54820250088277
$arrGI[2,2] = 19.94143858370069
$arrGI[3,0] = 33.1813899863129
$arrGI[3,1] = 14.54240228988833
$arrGI[3,2] = 19.94368203263307
$arrGI[4,0] = 33.14782424475341
$arrGI[4,1] = 14.54151014352986
$arrGI[4,2] = 19.94416165169323
$arrGI[5,0] = 33.38158494880805
$arrGI[5,1] = 14.54811951829107
$arrGI[5,2] = 19.94146165297894
$arrGI[6,0] = 34.42992415949547
$arrGI[6,1] = 14.58760925026978
$arrGI[6,2] = 19.94518885964003
$arrGI[7,0] = 36.52327189201864
$arrGI[7,1] = 14.7021199687355
$arrGI[7,2] = 20.00930739698332
$arrGI[8,0] = 38.06454537529351
$arrGI[8,1] = 14.80870158253534
$arrGI[8,2] = 20.0914416857313
$arrGI[9,0] = 38.76255524771378
$arrGI[9,1] = 14.86200178662173
$arrGI[9,2] = 20.1364586957641
$arrGI[10,0] = 39.02612607960653
$arrGI[10,1] = 14.88287021517906
$arrGI[10,2] = 20.15460676451123
$arrGI[11,0] = 38.96939934170216
$arrGI[11,1] = 14.87834550579973
$arrGI[11,2] = 20.15064929320517
$arrGI[12,0] = 38.78425607163746
$arrGI[12,1] = 14.86370496926238
$arrGI[12,2] = 20.13792968434311
$arrGI[13,0] = 38.67074415288572
$arrGI[13,1] = 14.85482614756234
$arrGI[13,2] = 20.13028193574656
$arrGI[14,0] = 38.01883974010272
$arrGI[14,1] = 14.80531456614461
$arrGI[14,2] = 20.08865381721393
$arrGI[15,0] = 37.61789701879588
$arrGI[15,1] = 14.77616918664281
$arrGI[15,2] = 20.06507757668678
$arrGI[16,0] = 37.38700821281438
$arrGI[16,1] = 14.75985935950125
$arrGI[16,2] = 20.05223764044489
$arrGI[17,0] = 37.30879512965004
$arrGI[17,1] = 14.75441528913094
$arrGI[17,2] = 20.0480139574866
$arrGI[18,0] = 37.66060898580791
$arrGI[18,1] = 14.77922484441289
$arrGI[18,2] = 20.06751272969804
$arrGI[19,0] = 38.83865980446527
$arrGI[19,1] = 14.86798673245299
$arrGI[19,2] = 20.14163586343305
$arrGI[20,0] = 39.60408498315896
$arrGI[20,1] = 14.92998478343873
$arrGI[20,2] = 20.19649659972555
$arrGI[21,0] = 39.19606916107776
$arrGI[21,1] = 14.89653332819168
$arrGI[21,2] = 20.1666295564735
$arrGI[22,0] = 37.64130007071385
$arrGI[22,1] = 14.77784199103057
$arrGI[22,2] = 20.06640957201555
$arrGI[23,0] = 35.95501888430339
$arrGI[23,1] = 14.66718272781677
$arrGI[23,2] = 19.98583199422967
$ws.Range("G2:I25").Value = $arrGI

$arrN = New-Object 'object[,]' 24,1
$arrN[0,0] = 16.21006386266351
$arrN[1,0] = 16.28718261907238
$arrN[2,0] = 16.33644568882588
$arrN[3,0] = 16.35700381192694
$arrN[4,0] = 16.36044671245668
$arrN[5,0] = 16.33672098421568
$arrN[6,0] = 16.23625898540869
$arrN[7,0] = 16.05431713006008
$arrN[8,0] = 15.92967672438214
$arrN[9,0] = 15.87490315991306
$arrN[10,0] = 15.85443629091401
$arrN[11,0] = 15.85883201645597
$arrN[12,0] = 15.87321384605146
$arrN[13,0] = 15.88205883756833
$arrN[14,0] = 15.93329486801004
$arrN[15,0] = 15.96521819491699
$arrN[16,0] = 15.98376106635576
$arrN[17,0] = 15.99007059154718
$arrN[18,0] = 15.96180114256252
$arrN[19,0] = 15.86898211709625
$arrN[20,0] = 15.8099195738547
$arrN[21,0] = 15.84129670864666
$arrN[22,0] = 15.96334540080818
$arrN[23,0] = 16.10194013523804
$ws.Range("N2:N25").Value = $arrN

Write-Host "Applied loading_percent updates"